$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column BO (col 67) holds "nomor_s" values for rows 6 through 196.
# Update them all from 2010 to 2338.
$ws.Range("BO6:BO196").Value = 2338
